$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 250000590
$ws.Range("I32").Value = 500000350
$ws.Range("J32").Value = 840
$ws.Range("K32").Value = 500000350
$ws.Range("L32").Value = 840
$ws.Range("M32").Value = -500000024
$ws.Range("N32").Value = -1492

$ws.Range("H121").Value = 1204.1666
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1204.1666
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("M121").Value = 3612.4998
$ws.Range("N121").Value = -7106.4998

$ws.Range("H137").Value = 3659877.8
$ws.Range("J137").Value = 8335242
$ws.Range("L137").Value = 25005726
$ws.Range("N137").Value = -25010826

$ws.Range("H138").Value = 2323.4727
$ws.Range("I138").Value = 2421.2354
$ws.Range("J138").Value = 2279.7368
$ws.Range("K138").Value = 7263.706200000001
$ws.Range("L138").Value = 6839.2104
$ws.Range("M138").Value = -2123.706200000001
$ws.Range("N138").Value = -17119.2104

$ws.Range("H141").Value = 2587.5938
$ws.Range("I141").Value = 1836.2632
$ws.Range("J141").Value = 3685.6924
$ws.Range("K141").Value = 5508.7896
$ws.Range("L141").Value = 11057.0772
$ws.Range("M141").Value = -328.7896000000001
$ws.Range("N141").Value = -21417.0772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3068.3333
$ws.Range("I2").Value = 1081.3889
$ws.Range("J2").Value = 7042.222
$ws.Range("K2").Value = 1081.3889
$ws.Range("L2").Value = 7042.222
$ws.Range("M2").Value = -968.3888999999999
$ws.Range("N2").Value = -7268.222

$ws.Range("H32").Value = 25437.951
$ws.Range("I32").Value = 5128.771
$ws.Range("K32").Value = 5128.771
$ws.Range("M32").Value = -4841.771

$ws.Range("H45").Value = 2149.5
$ws.Range("I45").Value = 1526.1
$ws.Range("K45").Value = 1526.1
$ws.Range("M45").Value = -1149.1

$ws.Range("H97").Value = 988.4516
$ws.Range("I97").Value = 814.5454999999999
$ws.Range("J97").Value = 1413.5555
$ws.Range("K97").Value = 814.5454999999999
$ws.Range("L97").Value = 1413.5555
$ws.Range("M97").Value = -318.5454999999999
$ws.Range("N97").Value = -2405.5555

$ws.Range("H110").Value = 3298.5642
$ws.Range("I110").Value = 2980.8064
$ws.Range("J110").Value = 4529.875
$ws.Range("K110").Value = 2980.8064
$ws.Range("L110").Value = 4529.875
$ws.Range("M110").Value = -935.8063999999999
$ws.Range("N110").Value = -8619.875

$ws.Range("H116").Value = 3068.3333
$ws.Range("I116").Value = 1081.3889
$ws.Range("J116").Value = 7042.222
$ws.Range("K116").Value = 1081.3889
$ws.Range("L116").Value = 7042.222
$ws.Range("M116").Value = 1212.6111
$ws.Range("N116").Value = -11630.222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3068.3333
$ws.Range("I3").Value = 1081.3889
$ws.Range("J3").Value = 7042.222
$ws.Range("K3").Value = 1081.3889
$ws.Range("L3").Value = 7042.222
$ws.Range("M3").Value = -967.3888999999999
$ws.Range("N3").Value = -7270.222

$ws.Range("H99").Value = 1636.0625
$ws.Range("I99").Value = 1032.9445
$ws.Range("J99").Value = 2411.5
$ws.Range("K99").Value = 1032.9445
$ws.Range("L99").Value = 2411.5
$ws.Range("M99").Value = 465.0554999999999
$ws.Range("N99").Value = -5407.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1334.2667
$ws.Range("I16").Value = 1276.1666
$ws.Range("J16").Value = 1566.6666
$ws.Range("K16").Value = 1276.1666
$ws.Range("L16").Value = 1566.6666
$ws.Range("M16").Value = -989.1666
$ws.Range("N16").Value = -2140.6666

$ws.Range("H31").Value = 1791.0952
$ws.Range("I31").Value = 1591.3024
$ws.Range("J31").Value = 2220.65
$ws.Range("K31").Value = 1591.3024
$ws.Range("L31").Value = 2220.65
$ws.Range("M31").Value = -1296.3024
$ws.Range("N31").Value = -2810.65

$ws.Range("H34").Value = 1791.0952
$ws.Range("I34").Value = 1591.3024
$ws.Range("J34").Value = 2220.65
$ws.Range("K34").Value = 1591.3024
$ws.Range("L34").Value = 2220.65
$ws.Range("M34").Value = -1389.3024
$ws.Range("N34").Value = -2624.65

$ws.Range("H113").Value = 1334.2667
$ws.Range("I113").Value = 1276.1666
$ws.Range("J113").Value = 1566.6666
$ws.Range("K113").Value = 1276.1666
$ws.Range("L113").Value = 1566.6666
$ws.Range("M113").Value = 893.8334
$ws.Range("N113").Value = -5906.6666

$ws.Range("H134").Value = 4262.6055
$ws.Range("I134").Value = 5012.2334
$ws.Range("J134").Value = 1451.5
$ws.Range("K134").Value = 15036.7002
$ws.Range("L134").Value = 4354.5
$ws.Range("M134").Value = -12501.7002
$ws.Range("N134").Value = -9424.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3745.5925
$ws.Range("I117").Value = 1000
$ws.Range("J117").Value = 3965.24
$ws.Range("K117").Value = 3000
$ws.Range("L117").Value = 11895.72
$ws.Range("M117").Value = 442
$ws.Range("N117").Value = -18779.72

$ws.Range("H129").Value = 1868.0869
$ws.Range("I129").Value = 913.3333
$ws.Range("J129").Value = 2205.0588
$ws.Range("K129").Value = 2739.9999
$ws.Range("L129").Value = 6615.176399999999
$ws.Range("M129").Value = 2260.0001
$ws.Range("N129").Value = -16615.1764

$ws.Range("H131").Value = 870.39703
$ws.Range("I131").Value = 615
$ws.Range("J131").Value = 886.3594000000001
$ws.Range("K131").Value = 1845
$ws.Range("L131").Value = 2659.0782
$ws.Range("M131").Value = 3195
$ws.Range("N131").Value = -12739.0782

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1260.0769
$ws.Range("I113").Value = 1159.95
$ws.Range("J113").Value = 1593.8334
$ws.Range("K113").Value = 1159.95
$ws.Range("L113").Value = 1593.8334
$ws.Range("M113").Value = 1010.05
$ws.Range("N113").Value = -5933.8334

$ws.Range("H126").Value = 4714.207
$ws.Range("I126").Value = 2600.75
$ws.Range("J126").Value = 7315.385
$ws.Range("K126").Value = 7802.25
$ws.Range("L126").Value = 21946.155
$ws.Range("M126").Value = -5332.25
$ws.Range("N126").Value = -26886.155

$ws.Range("H132").Value = 2393.0435
$ws.Range("I132").Value = 2092.111
$ws.Range("J132").Value = 3476.4
$ws.Range("K132").Value = 6276.333
$ws.Range("L132").Value = 10429.2
$ws.Range("M132").Value = -3746.333
$ws.Range("N132").Value = -15489.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1174.8214
$ws.Range("I46").Value = 998.2222
$ws.Range("J46").Value = 1492.7
$ws.Range("K46").Value = 998.2222
$ws.Range("L46").Value = 1492.7
$ws.Range("M46").Value = -810.2222
$ws.Range("N46").Value = -1868.7

$ws.Range("H55").Value = 739.7
$ws.Range("I55").Value = 418.375
$ws.Range("K55").Value = 418.375
$ws.Range("M55").Value = -245.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1859.381
$ws.Range("I107").Value = 919.1111
$ws.Range("J107").Value = 7501
$ws.Range("K107").Value = 2757.3333
$ws.Range("L107").Value = 22503
$ws.Range("M107").Value = -837.3332999999998
$ws.Range("N107").Value = -26343

$ws.Range("H122").Value = 3256.8857
$ws.Range("I122").Value = 2876.0476
$ws.Range("J122").Value = 3828.1428
$ws.Range("K122").Value = 8628.1428
$ws.Range("L122").Value = 11484.4284
$ws.Range("M122").Value = -6178.1428
$ws.Range("N122").Value = -16384.4284

$ws.Range("H132").Value = 2375.756
$ws.Range("I132").Value = 1934.7778
$ws.Range("J132").Value = 3226.2144
$ws.Range("K132").Value = 5804.3334
$ws.Range("L132").Value = 9678.643199999999
$ws.Range("M132").Value = -3274.3334
$ws.Range("N132").Value = -14738.6432
